# Add "RES boundary attack" columns (AS:AZ) for SEED 0
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header block (merged AS1:AZ1, labelled "BOUNDARY") ---
# Merge first, then paste the formatting, so the merge doesn't strip the
# inner borders back out into separate style records.
$ws.Range("AS1:AZ1").Merge()
$ws.Range("AR1").Copy()
$ws.Range("AS1:AZ1").PasteSpecial(-4122)
$ws.Range("AS1").Value = "BOUNDARY"

# --- Row 2: epsilon labels (stored as text, like the other epsilon rows) ---
$ws.Range("AR2").Copy()
$ws.Range("AS2:AZ2").PasteSpecial(-4122)
$ws.Range("AS2:AZ2").NumberFormat = "@"
$ws.Range("AS2").Value = "0.01"
$ws.Range("AT2").Value = "0.02"
$ws.Range("AU2").Value = "0.03"
$ws.Range("AV2").Value = "0.04"
$ws.Range("AW2").Value = "0.05"
$ws.Range("AX2").Value = "0.07"
$ws.Range("AY2").Value = "0.10"
$ws.Range("AZ2").Value = "0.20"

# --- Row 4: LSTM / MAE ---
$ws.Range("AS4").Value = 5.496322631835938
$ws.Range("AT4").Value = 5.549904823303223
$ws.Range("AU4").Value = 5.578744411468506
$ws.Range("AV4").Value = 5.667272090911865
$ws.Range("AW4").Value = 5.76094388961792
$ws.Range("AX4").Value = 5.969589710235596
$ws.Range("AY4").Value = 6.133481979370117
$ws.Range("AZ4").Value = 8.110892295837402

# --- Row 5: LSTM / RMSE ---
$ws.Range("AS5").Value = 6.761733136087383
$ws.Range("AT5").Value = 6.816339932968947
$ws.Range("AU5").Value = 6.823310044768669
$ws.Range("AV5").Value = 6.97236382221174
$ws.Range("AW5").Value = 7.045091356777554
$ws.Range("AX5").Value = 7.354597934544661
$ws.Range("AY5").Value = 7.627962881301
$ws.Range("AZ5").Value = 9.988832723336831

# --- Row 6: LSTM / SIM ---
$ws.Range("AS6").Value = 0.9996209740638733
$ws.Range("AT6").Value = 0.9996138215065002
$ws.Range("AU6").Value = 0.9996082782745361
$ws.Range("AV6").Value = 0.9995941519737244
$ws.Range("AW6").Value = 0.9995812177658081
$ws.Range("AX6").Value = 0.9995065927505493
$ws.Range("AY6").Value = 0.9994520545005798
$ws.Range("AZ6").Value = 0.99898761510849

# --- Row 7: RNN / MAE ---
$ws.Range("AS7").Value = 3.20345139503479
$ws.Range("AT7").Value = 3.456241846084595
$ws.Range("AU7").Value = 3.82329797744751
$ws.Range("AV7").Value = 3.98072338104248
$ws.Range("AW7").Value = 4.381083488464355
$ws.Range("AX7").Value = 5.53473424911499
$ws.Range("AY7").Value = 7.386642456054688
$ws.Range("AZ7").Value = 11.75538921356201

# --- Row 8: RNN / RMSE ---
$ws.Range("AS8").Value = 4.134818442993423
$ws.Range("AT8").Value = 4.465996167165476
$ws.Range("AU8").Value = 4.772953989168013
$ws.Range("AV8").Value = 5.088784073125476
$ws.Range("AW8").Value = 5.650556027219626
$ws.Range("AX8").Value = 6.979809388889112
$ws.Range("AY8").Value = 9.280631407599941
$ws.Range("AZ8").Value = 14.79750866927885

# --- Row 9: RNN / SIM ---
$ws.Range("AS9").Value = 0.9997895359992981
$ws.Range("AT9").Value = 0.9997361898422241
$ws.Range("AU9").Value = 0.9996867775917053
$ws.Range("AV9").Value = 0.9996429085731506
$ws.Range("AW9").Value = 0.9994867444038391
$ws.Range("AX9").Value = 0.9991776347160339
$ws.Range("AY9").Value = 0.9985546469688416
$ws.Range("AZ9").Value = 0.9961950778961182

# --- Row 10: GRU / MAE ---
$ws.Range("AS10").Value = 2.61094069480896
$ws.Range("AT10").Value = 2.756315469741821
$ws.Range("AU10").Value = 3.124436855316162
$ws.Range("AV10").Value = 3.192680358886719
$ws.Range("AW10").Value = 3.359122514724731
$ws.Range("AX10").Value = 4.486050128936768
$ws.Range("AY10").Value = 5.294600963592529
$ws.Range("AZ10").Value = 9.870691299438477

# --- Row 11: GRU / RMSE ---
$ws.Range("AS11").Value = 3.470677433122761
$ws.Range("AT11").Value = 3.673629795949171
$ws.Range("AU11").Value = 4.042502615036646
$ws.Range("AV11").Value = 4.205121895455546
$ws.Range("AW11").Value = 4.298416804350851
$ws.Range("AX11").Value = 5.758591908428854
$ws.Range("AY11").Value = 6.569731878328446
$ws.Range("AZ11").Value = 12.43660437651567

# --- Row 12: GRU / SIM ---
$ws.Range("AS12").Value = 0.9997884631156921
$ws.Range("AT12").Value = 0.9997621774673462
$ws.Range("AU12").Value = 0.9997113943099976
$ws.Range("AV12").Value = 0.9996863603591919
$ws.Range("AW12").Value = 0.9996731281280518
$ws.Range("AX12").Value = 0.9994210600852966
$ws.Range("AY12").Value = 0.9992382526397705
$ws.Range("AZ12").Value = 0.9972619414329529
